# Auto-generated Excel COM-interop script to apply market-data refresh edits
# across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets (columns H-N, per-row).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# ALC row 12
$ws.Range("H12").Value = 1338.3846
$ws.Range("I12").Value = 465.16666
$ws.Range("J12").Value = 2086.8572
$ws.Range("K12").Value = 465.16666
$ws.Range("L12").Value = 2086.8572
$ws.Range("M12").Value = -295.16666

# ALC row 86
$ws.Range("H86").Value = 3264
$ws.Range("I86").Value = 1894
$ws.Range("J86").Value = 6004
$ws.Range("K86").Value = 1894
$ws.Range("L86").Value = 6004
$ws.Range("M86").Value = -771
$ws.Range("N86").Value = -8250

# ALC row 89
$ws.Range("H89").Value = 3264
$ws.Range("I89").Value = 1894
$ws.Range("J89").Value = 6004
$ws.Range("K89").Value = 9470
$ws.Range("L89").Value = 30020
$ws.Range("M89").Value = -3854
$ws.Range("N89").Value = -41252

# ALC row 132
$ws.Range("H132").Value = 237318.19
$ws.Range("I132").Value = 246787.34
$ws.Range("J132").Value = 6902
$ws.Range("K132").Value = 740362.02
$ws.Range("L132").Value = 20706
$ws.Range("M132").Value = -737832.02

# ALC row 137
$ws.Range("H137").Value = 2088.4375
$ws.Range("I137").Value = 1176.25
$ws.Range("J137").Value = 2740
$ws.Range("K137").Value = 3528.75
$ws.Range("L137").Value = 8220
$ws.Range("M137").Value = -978.75

# ALC row 138
$ws.Range("H138").Value = 3240.255
$ws.Range("I138").Value = 2235.8125
$ws.Range("J138").Value = 3699.4285
$ws.Range("K138").Value = 6707.4375
$ws.Range("L138").Value = 11098.2855
$ws.Range("M138").Value = -1567.4375
$ws.Range("N138").Value = -21378.2855

$ws = $wb.Worksheets.Item("ARM")
# ARM row 4
$ws.Range("H4").Value = 1159.4
$ws.Range("I4").Value = 266.66666
$ws.Range("J4").Value = 2498.5
$ws.Range("K4").Value = 266.66666
$ws.Range("L4").Value = 2498.5
$ws.Range("M4").Value = -150.66666

# ARM row 48
$ws.Range("H48").Value = 0
$ws.Range("I48").Value = 0
$ws.Range("J48").Value = 0
$ws.Range("K48").Value = 0
$ws.Range("L48").Value = 0
$ws.Range("N48").Value = $null

# ARM row 63
$ws.Range("H63").Value = 2567
$ws.Range("I63").Value = 2567
$ws.Range("J63").Value = 0
$ws.Range("K63").Value = 2567
$ws.Range("L63").Value = 0
$ws.Range("M63").Value = -1881
$ws.Range("N63").Value = $null

# ARM row 66
$ws.Range("H66").Value = 2567
$ws.Range("I66").Value = 2567
$ws.Range("J66").Value = 0
$ws.Range("K66").Value = 12835
$ws.Range("L66").Value = 0
$ws.Range("M66").Value = -9403
$ws.Range("N66").Value = $null

# ARM row 74
$ws.Range("H74").Value = 3375.5974
$ws.Range("I74").Value = 2874.4407
$ws.Range("J74").Value = 5018.278
$ws.Range("K74").Value = 2874.4407
$ws.Range("L74").Value = 5018.278
$ws.Range("M74").Value = -2000.4407
$ws.Range("N74").Value = -6766.278

# ARM row 77
$ws.Range("H77").Value = 3375.5974
$ws.Range("I77").Value = 2874.4407
$ws.Range("J77").Value = 5018.278
$ws.Range("K77").Value = 14372.2035
$ws.Range("L77").Value = 25091.39
$ws.Range("M77").Value = -10004.2035
$ws.Range("N77").Value = -33827.39

# ARM row 97
$ws.Range("H97").Value = 911
$ws.Range("I97").Value = 815.375
$ws.Range("J97").Value = 1020.2857
$ws.Range("K97").Value = 815.375
$ws.Range("L97").Value = 1020.2857
$ws.Range("M97").Value = -319.375
$ws.Range("N97").Value = -2012.2857

# ARM row 102
$ws.Range("H102").Value = 3618.4375
$ws.Range("I102").Value = 3626.4
$ws.Range("J102").Value = 3499
$ws.Range("K102").Value = 3626.4
$ws.Range("L102").Value = 3499
$ws.Range("M102").Value = -2004.4

# ARM row 132
$ws.Range("H132").Value = 4811.2793
$ws.Range("I132").Value = 3122.225
$ws.Range("J132").Value = 7224.2144
$ws.Range("K132").Value = 9366.674999999999
$ws.Range("L132").Value = 21672.6432
$ws.Range("M132").Value = -6836.674999999999

$ws = $wb.Worksheets.Item("BSM")
# BSM row 61
$ws.Range("H61").Value = 20000
$ws.Range("I61").Value = 10000
$ws.Range("J61").Value = 30000
$ws.Range("K61").Value = 10000
$ws.Range("L61").Value = 30000
$ws.Range("M61").Value = -9687
$ws.Range("N61").Value = -30626

# BSM row 94
$ws.Range("H94").Value = 1160
$ws.Range("I94").Value = 729.8570999999999
$ws.Range("J94").Value = 2020.2858
$ws.Range("K94").Value = 729.8570999999999
$ws.Range("L94").Value = 2020.2858
$ws.Range("M94").Value = -278.8570999999999
$ws.Range("N94").Value = -2922.2858

# BSM row 107
$ws.Range("H107").Value = 1844
$ws.Range("I107").Value = 1882.2222
$ws.Range("J107").Value = 1500
$ws.Range("K107").Value = 1882.2222
$ws.Range("L107").Value = 1500
$ws.Range("M107").Value = 37.77780000000007
$ws.Range("N107").Value = -5340

# BSM row 119
$ws.Range("H119").Value = 70833.336
$ws.Range("I119").Value = 75000
$ws.Range("J119").Value = 70000
$ws.Range("K119").Value = 75000
$ws.Range("L119").Value = 70000
$ws.Range("M119").Value = -70162
$ws.Range("N119").Value = -79676

# BSM row 120
$ws.Range("H120").Value = 50000
$ws.Range("I120").Value = 0
$ws.Range("J120").Value = 50000
$ws.Range("K120").Value = 0
$ws.Range("L120").Value = 50000
$ws.Range("N120").Value = -59676

# BSM row 125
$ws.Range("H125").Value = 0
$ws.Range("I125").Value = 0
$ws.Range("J125").Value = 0
$ws.Range("K125").Value = 0
$ws.Range("L125").Value = 0
$ws.Range("N125").Value = $null

# BSM row 140
$ws.Range("H140").Value = 108714
$ws.Range("I140").Value = 0
$ws.Range("J140").Value = 108714
$ws.Range("K140").Value = 0
$ws.Range("L140").Value = 108714
$ws.Range("N140").Value = -119074

$ws = $wb.Worksheets.Item("CRP")
# CRP row 31
$ws.Range("H31").Value = 8437.049000000001
$ws.Range("I31").Value = 2186.5715
$ws.Range("J31").Value = 9723.912
$ws.Range("K31").Value = 2186.5715
$ws.Range("L31").Value = 9723.912
$ws.Range("M31").Value = -1891.5715
$ws.Range("N31").Value = -10313.912

# CRP row 34
$ws.Range("H34").Value = 8437.049000000001
$ws.Range("I34").Value = 2186.5715
$ws.Range("J34").Value = 9723.912
$ws.Range("K34").Value = 2186.5715
$ws.Range("L34").Value = 9723.912
$ws.Range("M34").Value = -1984.5715
$ws.Range("N34").Value = -10127.912

# CRP row 50
$ws.Range("H50").Value = 87851
$ws.Range("I50").Value = 0
$ws.Range("J50").Value = 87851
$ws.Range("K50").Value = 0
$ws.Range("L50").Value = 87851
$ws.Range("N50").Value = -89101

# CRP row 51
$ws.Range("H51").Value = 52499.75
$ws.Range("I51").Value = 0
$ws.Range("J51").Value = 52499.75
$ws.Range("K51").Value = 0
$ws.Range("L51").Value = 52499.75
$ws.Range("N51").Value = -53971.75

# CRP row 58
$ws.Range("H58").Value = 3267.6924
$ws.Range("I58").Value = 2874.524
$ws.Range("J58").Value = 4919
$ws.Range("K58").Value = 2874.524
$ws.Range("L58").Value = 4919
$ws.Range("M58").Value = -2671.524

# CRP row 61
$ws.Range("H61").Value = 52499.75
$ws.Range("I61").Value = 0
$ws.Range("J61").Value = 52499.75
$ws.Range("K61").Value = 0
$ws.Range("L61").Value = 52499.75
$ws.Range("N61").Value = -53195.75

# CRP row 107
$ws.Range("H107").Value = 1029.6428
$ws.Range("I107").Value = 1034.6666
$ws.Range("J107").Value = 999.5
$ws.Range("K107").Value = 1034.6666
$ws.Range("L107").Value = 999.5
$ws.Range("M107").Value = 885.3334
$ws.Range("N107").Value = -4839.5

# CRP row 134
$ws.Range("H134").Value = 22827.234
$ws.Range("I134").Value = 22827.234
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 68481.702
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = -65946.702

# CRP row 136
$ws.Range("H136").Value = 3267.6924
$ws.Range("I136").Value = 2874.524
$ws.Range("J136").Value = 4919
$ws.Range("K136").Value = 8623.572
$ws.Range("L136").Value = 14757
$ws.Range("M136").Value = -6073.572

$ws = $wb.Worksheets.Item("CUL")
# CUL row 92
$ws.Range("H92").Value = 858.35297
$ws.Range("I92").Value = 672.6
$ws.Range("J92").Value = 2251.5
$ws.Range("K92").Value = 2017.8
$ws.Range("L92").Value = 6754.5
$ws.Range("M92").Value = -769.8000000000002
$ws.Range("N92").Value = -9250.5

# CUL row 122
$ws.Range("H122").Value = 824.5454999999999
$ws.Range("I122").Value = 819
$ws.Range("J122").Value = 825.7778
$ws.Range("K122").Value = 7371
$ws.Range("L122").Value = 7432.000199999999
$ws.Range("M122").Value = -4921
$ws.Range("N122").Value = -12332.0002

$ws = $wb.Worksheets.Item("GSM")
# GSM row 80
$ws.Range("H80").Value = 458483
$ws.Range("I80").Value = 559190.5600000001
$ws.Range("J80").Value = 5299
$ws.Range("K80").Value = 559190.5600000001
$ws.Range("L80").Value = 5299
$ws.Range("M80").Value = -558192.5600000001
$ws.Range("N80").Value = -7295

# GSM row 83
$ws.Range("H83").Value = 458483
$ws.Range("I83").Value = 559190.5600000001
$ws.Range("J83").Value = 5299
$ws.Range("K83").Value = 2795952.8
$ws.Range("L83").Value = 26495
$ws.Range("M83").Value = -2790960.8
$ws.Range("N83").Value = -36479

# GSM row 107
$ws.Range("H107").Value = 463.53845
$ws.Range("I107").Value = 302.9
$ws.Range("J107").Value = 999
$ws.Range("K107").Value = 302.9
$ws.Range("L107").Value = 999
$ws.Range("M107").Value = 1617.1
$ws.Range("N107").Value = -4839

# GSM row 122
$ws.Range("H122").Value = 3753.6775
$ws.Range("I122").Value = 1908.1818
$ws.Range("J122").Value = 8264.888999999999
$ws.Range("K122").Value = 5724.5454
$ws.Range("L122").Value = 24794.667
$ws.Range("M122").Value = -3274.5454
$ws.Range("N122").Value = -29694.667

# GSM row 126
$ws.Range("H126").Value = 2264.3438
$ws.Range("I126").Value = 2096.318
$ws.Range("J126").Value = 2634
$ws.Range("K126").Value = 6288.954000000001
$ws.Range("L126").Value = 7902
$ws.Range("M126").Value = -3818.954000000001

$ws = $wb.Worksheets.Item("LTW")
# LTW row 7
$ws.Range("H7").Value = 3361.25
$ws.Range("I7").Value = 3280.8333
$ws.Range("J7").Value = 3602.5
$ws.Range("K7").Value = 3280.8333
$ws.Range("L7").Value = 3602.5
$ws.Range("M7").Value = -3168.8333

# LTW row 40
$ws.Range("H40").Value = 3788.6667
$ws.Range("I40").Value = 3683.5
$ws.Range("J40").Value = 3999
$ws.Range("K40").Value = 3683.5
$ws.Range("L40").Value = 3999
$ws.Range("M40").Value = -3547.5
$ws.Range("N40").Value = -4271

# LTW row 93
$ws.Range("H93").Value = 1160.2759
$ws.Range("I93").Value = 1073
$ws.Range("J93").Value = 1916.6666
$ws.Range("K93").Value = 1073
$ws.Range("L93").Value = 1916.6666
$ws.Range("M93").Value = 175

# LTW row 122
$ws.Range("H122").Value = 3886.1177
$ws.Range("I122").Value = 3654.5715
$ws.Range("J122").Value = 4966.6665
$ws.Range("K122").Value = 10963.7145
$ws.Range("L122").Value = 14899.9995
$ws.Range("M122").Value = -8513.7145
$ws.Range("N122").Value = -19799.9995

# LTW row 126
$ws.Range("H126").Value = 3361.25
$ws.Range("I126").Value = 3280.8333
$ws.Range("J126").Value = 3602.5
$ws.Range("K126").Value = 9842.499899999999
$ws.Range("L126").Value = 10807.5
$ws.Range("M126").Value = -7372.499899999999

# LTW row 136
$ws.Range("H136").Value = 6077.174
$ws.Range("I136").Value = 6533.6665
$ws.Range("J136").Value = 5579.1816
$ws.Range("K136").Value = 19600.9995
$ws.Range("L136").Value = 16737.5448
$ws.Range("M136").Value = -17050.9995
$ws.Range("N136").Value = -21837.5448

$ws = $wb.Worksheets.Item("WVR")
# WVR row 8
$ws.Range("H8").Value = 3333
$ws.Range("I8").Value = 0
$ws.Range("J8").Value = 3333
$ws.Range("K8").Value = 0
$ws.Range("L8").Value = 3333
$ws.Range("N8").Value = -3613

# WVR row 62
$ws.Range("H62").Value = 0
$ws.Range("I62").Value = 0
$ws.Range("J62").Value = 0
$ws.Range("K62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("M62").Value = $null

# WVR row 65
$ws.Range("H65").Value = 0
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("M65").Value = $null

# WVR row 122
$ws.Range("H122").Value = 3651.0476
$ws.Range("I122").Value = 3326.389
$ws.Range("J122").Value = 5599
$ws.Range("K122").Value = 9979.167000000001
$ws.Range("L122").Value = 16797
$ws.Range("M122").Value = -7529.167000000001

# WVR row 132
$ws.Range("H132").Value = 12152607
$ws.Range("I132").Value = 3957.2173
$ws.Range("J132").Value = 25458272
$ws.Range("K132").Value = 11871.6519
$ws.Range("L132").Value = 76374816
$ws.Range("M132").Value = -9341.651899999999
$ws.Range("N132").Value = -76379876

# WVR row 136
$ws.Range("H136").Value = 8060821.5
$ws.Range("I136").Value = 11180188
$ws.Range("J136").Value = 2457.6667
$ws.Range("K136").Value = 33540564
$ws.Range("L136").Value = 7373.000100000001
$ws.Range("M136").Value = -33538014
